# Apply updated crypto price/volume figures (scraper refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.869.24"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "2.461.02"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.460.57"
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +1.59%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.356"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "2.911.74"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "62.975.84"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("D18").Value = "2.457.92"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +1.00%  "
$ws.Range("E20").Value = "  +6.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +12.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.998"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "622.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000103"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.25%  "
$ws.Range("D29").Value = "2.567.92"
$ws.Range("E29").Value = "  +0.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("E31").Value = "  +5.79%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.142"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.09"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.45%  "
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "145.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.88%  "
$ws.Range("E42").Value = "  -1.01%  "
$ws.Range("E43").Value = "  +14.63%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("E46").Value = "  +2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "20.81"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.87%  "
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("E51").Value = "  -0.43%  "
